$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1087  # was 1086
$ws.Range("F8").Value = 259  # was 258
$ws.Range("F9").Value = 49  # was 48
$ws.Range("F10").Value = 1221  # was 1219
$ws.Range("F11").Value = 27642  # was 27630
$ws.Range("F12").Value = 3115  # was 3107
$ws.Range("F17").Value = 285  # was 284
$ws.Range("F18").Value = 551  # was 550
$ws.Range("F20").Value = 225  # was 223
$ws.Range("F26").Value = 474  # was 473
$ws.Range("F29").Value = 566  # was 565
$ws.Range("F31").Value = 27  # was 26

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 354  # was 353
$ws.Range("F7").Value = 679  # was 677
$ws.Range("F8").Value = 74  # was 72
$ws.Range("F12").Value = 177  # was 176
$ws.Range("F21").Value = 4207  # was 4205

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1096  # was 1095
$ws.Range("F5").Value = 286  # was 285

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1096  # was 1095
$ws.Range("F9").Value = 354  # was 353
$ws.Range("F10").Value = 286  # was 285
$ws.Range("F11").Value = 679  # was 677
$ws.Range("F12").Value = 1087  # was 1086
$ws.Range("F16").Value = 259  # was 258
$ws.Range("F17").Value = 49  # was 48
$ws.Range("F18").Value = 1221  # was 1219
$ws.Range("F19").Value = 27642  # was 27630
$ws.Range("F20").Value = 74  # was 72
$ws.Range("F23").Value = 177  # was 176
$ws.Range("F26").Value = 3115  # was 3107
$ws.Range("F33").Value = 285  # was 284
$ws.Range("F34").Value = 551  # was 550
$ws.Range("F46").Value = 566  # was 565
$ws.Range("F48").Value = 4207  # was 4205
$ws.Range("F49").Value = 27  # was 26

